$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = -0.116
$ws.Range("G2").Value = -0.03794392523364486
$ws.Range("H2").Value = -0.03794392523364486
$ws.Range("I2").Value = -0.2453133194285434
$ws.Range("J2").Value = -0.2453133194285434
$ws.Range("L2").Value = -0.1074766355140187
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 37.3
$ws.Range("V2").Value = 1.643171806167401
$ws.Range("W2").Value = -0.02725118483412322
$ws.Range("X2").Value = 0.06356401794851509
$ws.Range("Y2").Value = -0.0908152027826383
$ws.Range("Z2").Value = 0.9766104009158767
$ws.Range("AA2").Value = -0.2395755392371143
$ws.Range("AB2").Value = 0.06347733752046189
$ws.Range("AC2").Value = -0.3030528767575762
$ws.Range("AE2").Value = 0.05426258942707171
$ws.Range("AF2").Value = 0.05426258942707171
$ws.Range("AG2").Value = -37.24573741057293
$ws.Range("AH2").Value = 0.002384721948857407
$ws.Range("AI2").Value = 0.001331458011490291
$ws.Range("AJ2").Value = 2.560594651151886
$ws.Range("AK2").Value = -10.78254372570747
$ws.Range("AP2").Value = 18.96422475080088
$ws.Range("AQ2").Value = 2.495327102803738

# Row 3 updates
$ws.Range("B3").Value = "Automated Systems Company - KPSC (KWSE:ASC)"
$ws.Range("D3").Value = -0.116
$ws.Range("G3").Value = -0.03794392523364486
$ws.Range("H3").Value = -0.03794392523364486
$ws.Range("I3").Value = -0.2453133194285434
$ws.Range("J3").Value = -0.2453133194285434
$ws.Range("L3").Value = -0.1074766355140187
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 37.3
$ws.Range("V3").Value = 1.643171806167401
$ws.Range("W3").Value = -0.02725118483412322
$ws.Range("X3").Value = 0.06356401794851509
$ws.Range("Y3").Value = -0.0908152027826383
$ws.Range("Z3").Value = 0.9766104009158767
$ws.Range("AA3").Value = -0.2395755392371143
$ws.Range("AB3").Value = 0.06347733752046189
$ws.Range("AC3").Value = -0.3030528767575762
$ws.Range("AE3").Value = 0.05426258942707171
$ws.Range("AF3").Value = 0.05426258942707171
$ws.Range("AG3").Value = -37.24573741057293
$ws.Range("AH3").Value = 0.002384721948857407
$ws.Range("AI3").Value = 0.001331458011490291
$ws.Range("AJ3").Value = 2.560594651151886
$ws.Range("AK3").Value = -10.78254372570747
$ws.Range("AP3").Value = 18.96422475080088
$ws.Range("AQ3").Value = 2.495327102803738
